# Auto-generated edit script applying the commit diff to Sagittarius_Profits workbook.
# For each sheet/cell: either set a new numeric value, or clear the cell entirely
# (ClearContents) when the cell was removed in the diff, or set a value on a
# previously-empty cell when the cell was newly added in the diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 1000
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H29").Value = 7
$ws.Range("I29").Value = 7
$ws.Range("K29").Value = 21
$ws.Range("M29").Value = 260
$ws.Range("H34").Value = 9492
$ws.Range("I34").Value = 9492
$ws.Range("K34").Value = 9492
$ws.Range("M34").Value = -9289
$ws.Range("H36").Value = 9492
$ws.Range("I36").Value = 9492
$ws.Range("K36").Value = 9492
$ws.Range("M36").Value = -8777
$ws.Range("H39").Value = 523.9231
$ws.Range("I39").Value = 615.63635
$ws.Range("J39").Value = 19.5
$ws.Range("K39").Value = 1846.90905
$ws.Range("L39").Value = 58.5
$ws.Range("M39").Value = -1550.90905
$ws.Range("N39").Value = -650.5
$ws.Range("H47").Value = 17632.666
$ws.Range("J47").Value = 17949
$ws.Range("L47").Value = 17949
$ws.Range("N47").Value = -19893
$ws.Range("H54").Value = 62499
$ws.Range("J54").Value = 24999
$ws.Range("L54").Value = 24999
$ws.Range("N54").Value = -25971
$ws.Range("H98").Value = 2047.0625
$ws.Range("I98").Value = 593.5714
$ws.Range("K98").Value = 593.5714
$ws.Range("M98").Value = 904.4286
$ws.Range("H113").Value = 4275
$ws.Range("J113").Value = 4500
$ws.Range("L113").Value = 4500
$ws.Range("N113").Value = -11008
$ws.Range("H122").Value = 2047.0625
$ws.Range("I122").Value = 593.5714
$ws.Range("K122").Value = 1780.7142
$ws.Range("M122").Value = 669.2857999999999
$ws.Range("H129").Value = 4142.1113
$ws.Range("I129").Value = 4998.25
$ws.Range("J129").Value = 3457.2
$ws.Range("K129").Value = 14994.75
$ws.Range("L129").Value = 10371.6
$ws.Range("M129").Value = -9994.75
$ws.Range("N129").Value = -20371.6
$ws.Range("H137").Value = 1802.88
$ws.Range("J137").Value = 1799
$ws.Range("L137").Value = 5397
$ws.Range("N137").Value = -10497
$ws.Range("H138").Value = 2052.2896
$ws.Range("I138").Value = 1797.3448
$ws.Range("J138").Value = 2873.7778
$ws.Range("K138").Value = 5392.0344
$ws.Range("L138").Value = 8621.3334
$ws.Range("M138").Value = -252.0344000000005
$ws.Range("N138").Value = -18901.3334
$ws.Range("H141").Value = 1410.0869
$ws.Range("I141").Value = 1442.8636
$ws.Range("K141").Value = 4328.5908
$ws.Range("M141").Value = 851.4092000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1088.5834
$ws.Range("I5").Value = 1406.7778
$ws.Range("K5").Value = 1406.7778
$ws.Range("M5").Value = -1294.7778
$ws.Range("H32").Value = 7224.2856
$ws.Range("I32").Value = 6428.4165
$ws.Range("J32").Value = 11999.5
$ws.Range("K32").Value = 6428.4165
$ws.Range("L32").Value = 11999.5
$ws.Range("M32").Value = -6141.4165
$ws.Range("N32").Value = -12573.5
$ws.Range("H43").Value = 42499
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H61").Value = 2215
$ws.Range("I61").Value = 2215
$ws.Range("K61").Value = 2215
$ws.Range("M61").Value = -2003
$ws.Range("H74").Value = 2707.8572
$ws.Range("I74").Value = 2742.5
$ws.Range("K74").Value = 2742.5
$ws.Range("M74").Value = -1868.5
$ws.Range("H77").Value = 2707.8572
$ws.Range("I77").Value = 2742.5
$ws.Range("K77").Value = 13712.5
$ws.Range("M77").Value = -9344.5
$ws.Range("H97").Value = 342.31818
$ws.Range("I97").Value = 374.7647
$ws.Range("K97").Value = 374.7647
$ws.Range("M97").Value = 121.2353
$ws.Range("H102").Value = 2125.3333
$ws.Range("I102").Value = 2125.3333
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2125.3333
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -503.3332999999998
$ws.Range("H136").Value = 2215
$ws.Range("I136").Value = 2215
$ws.Range("K136").Value = 6645
$ws.Range("M136").Value = -4095

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1088.5834
$ws.Range("I4").Value = 1406.7778
$ws.Range("K4").Value = 1406.7778
$ws.Range("M4").Value = -1291.7778
$ws.Range("H105").Value = 3750
$ws.Range("I105").Value = 3750
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 3750
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -2003
$ws.Range("H134").Value = 1793
$ws.Range("I134").Value = 1906.4286
$ws.Range("J134").Value = 999
$ws.Range("K134").Value = 5719.2858
$ws.Range("L134").Value = 2997
$ws.Range("M134").Value = -3184.2858
$ws.Range("N134").Value = -8067

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2214.3333
$ws.Range("I31").Value = 1863.2
$ws.Range("K31").Value = 1863.2
$ws.Range("M31").Value = -1568.2
$ws.Range("H34").Value = 2214.3333
$ws.Range("I34").Value = 1863.2
$ws.Range("K34").Value = 1863.2
$ws.Range("M34").Value = -1661.2
$ws.Range("H62").Value = 1315
$ws.Range("I62").Value = 1315
$ws.Range("K62").Value = 1315
$ws.Range("M62").Value = -691
$ws.Range("H65").Value = 1315
$ws.Range("I65").Value = 1315
$ws.Range("K65").Value = 6575
$ws.Range("M65").Value = -3455
$ws.Range("H132").Value = 3344.111
$ws.Range("I132").Value = 3399.5715
$ws.Range("K132").Value = 10198.7145
$ws.Range("M132").Value = -7668.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 8462029
$ws.Range("I4").Value = 9166782
$ws.Range("K4").Value = 27500346
$ws.Range("M4").Value = -27500234
$ws.Range("H140").Value = 9835.75
$ws.Range("I140").Value = 2156.2727
$ws.Range("K140").Value = 6468.8181
$ws.Range("M140").Value = -1288.8181
$ws.Range("H141").Value = 6032.8335
$ws.Range("I141").Value = 6032.8335
$ws.Range("K141").Value = 18098.5005
$ws.Range("M141").Value = -12918.5005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6960.1113
$ws.Range("I70").Value = 6600
$ws.Range("K70").Value = 6600
$ws.Range("M70").Value = -6330
$ws.Range("H73").Value = 6960.1113
$ws.Range("I73").Value = 6600
$ws.Range("K73").Value = 6600
$ws.Range("M73").Value = -5664
$ws.Range("H122").Value = 3147.25
$ws.Range("I122").Value = 2195
$ws.Range("J122").Value = 4099.5
$ws.Range("K122").Value = 6585
$ws.Range("L122").Value = 12298.5
$ws.Range("M122").Value = -4135
$ws.Range("N122").Value = -17198.5
$ws.Range("H132").Value = 2737.25
$ws.Range("I132").Value = 2737.25
$ws.Range("K132").Value = 8211.75
$ws.Range("M132").Value = -5681.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1523.1538
$ws.Range("I22").Value = 1400.0834
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 1400.0834
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -1105.0834
$ws.Range("N22").Value = -3590
$ws.Range("H27").Value = 1523.1538
$ws.Range("I27").Value = 1400.0834
$ws.Range("J27").Value = 3000
$ws.Range("K27").Value = 1400.0834
$ws.Range("L27").Value = 3000
$ws.Range("M27").Value = -1293.0834
$ws.Range("N27").Value = -3214
$ws.Range("H122").Value = 6188
$ws.Range("I122").Value = 9293.75
$ws.Range("J122").Value = 4726.4707
$ws.Range("K122").Value = 27881.25
$ws.Range("L122").Value = 14179.4121
$ws.Range("M122").Value = -25431.25
$ws.Range("N122").Value = -19079.4121
$ws.Range("H132").Value = 1998.7693
$ws.Range("I132").Value = 1998.7693
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5996.3079
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -3466.3079
$ws.Range("H136").Value = 4060.3333
$ws.Range("I136").Value = 3148.4285
$ws.Range("K136").Value = 9445.2855
$ws.Range("M136").Value = -6895.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 100626
$ws.Range("J108").Value = 100626
$ws.Range("L108").Value = 100626
$ws.Range("N108").Value = -108306
$ws.Range("H113").Value = 1859
$ws.Range("I113").Value = 1873.5
$ws.Range("K113").Value = 5620.5
$ws.Range("M113").Value = -3450.5
$ws.Range("H132").Value = 2984.652
$ws.Range("I132").Value = 3134.3684
$ws.Range("K132").Value = 9403.1052
$ws.Range("M132").Value = -6873.1052
$ws.Range("H136").Value = 3271.7368
$ws.Range("I136").Value = 3259.0557
$ws.Range("K136").Value = 9777.167099999999
$ws.Range("M136").Value = -7227.167099999999
